# (JMT) Updated 1s6 to meet acceleration specification:
# insert a new "npc_speed_increase" parameter column (with value
# "[-5..-30]kph") right after "npc_speed_gt_dut_start", pushing the
# remaining parameter columns one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (6th column: was
# npc_dist_lt_dut_start), shifting the old F:J columns to G:K.
$ws.Columns.Item(6).Insert() | Out-Null

# Fill in the new column's header/value. Value first, then header, so the
# shared-string table picks up the two new strings in the same order as
# the authored workbook ("[-5..-30]kph" before "npc_speed_increase").
$ws.Cells.Item(3, 6).Value = "[-5..-30]kph"
$ws.Cells.Item(2, 6).Value = "npc_speed_increase"

# Give the new column the same width as column E (its neighbour); the
# shifted G:K columns keep their original best-fit widths automatically.
$ws.Columns.Item(6).ColumnWidth = 22.28515625

# Match the saved selection/view state.
$ws.Range("F10").Select() | Out-Null
